# Regenerate merged AHB files:
#  1) rename the "_old" / "_new" header-label suffixes to "_FV2310" / "_FV2404"
#  2) freeze the header row
#  3) turn the data range into a proper Excel Table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row relabelling -------------------------------------------
# Columns A:J were "<Label>_old", columns L:U were "<Label>_new" (column K
# holds the constant "diff" header and is left untouched).
$fv2310Labels = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $fv2310Labels.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($fv2310Labels[$i])_FV2310"
    $ws.Cells.Item(1, $i + 12).Value = "$($fv2310Labels[$i])_FV2404"
}

# --- 2. Freeze panes at the header row ------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn A1:U92 into an Excel Table -----------------------------------
$tableRange = $ws.Range("A1:U92")
$table = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$table.Name = "Table1"
$table.TableStyle = ""

Write-Host "Done"
